# "Working on 2023 report"
# Append water-year rows for WY 2019-2023 to the Yr-type table on Sheet1.
# Only columns A (WY), E (Index) and F (Yr-type) have data for these new
# rows -- B/C/D (Oct-Mar / Apr-Jul / WYsum) are left blank, same as the
# source data feeding this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 115; WY = 2019; Index = 10.34; YrType = "W" },
    @{ Row = 116; WY = 2020; Index = 6.13;  YrType = "D" },
    @{ Row = 117; WY = 2021; Index = 3.86;  YrType = "C" },
    @{ Row = 118; WY = 2022; Index = 4.5;   YrType = "C" },
    @{ Row = 119; WY = 2023; Index = 9.35;  YrType = "W" }
)

# The last existing data row (114) carries the formatting (style index 1:
# left-aligned Consolas) used for every WY cell in column A. Copy that
# formatting onto each new WY cell instead of re-deriving it, so we don't
# fork a brand-new style entry in styles.xml.
$formatSource = $ws.Range("A114")

foreach ($r in $newRows) {
    $rowNum = $r.Row

    [void]$formatSource.Copy()
    $aCell = $ws.Range("A$rowNum")
    [void]$aCell.PasteSpecial(-4122)  # xlPasteFormats
    $aCell.Value2 = $r.WY

    $ws.Range("E$rowNum").Value2 = $r.Index
    $ws.Range("F$rowNum").Value2 = $r.YrType
}

$excel.CutCopyMode = $false

# Match the author's final selection / scroll position after entering the
# 2023 row.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 95
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F120").Select()
